$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M")

$data = @{
    2 = @(1.02,1.029170723896049,1.037835783583608,1.029037989954252,1.044729243049703,1.033768717741986,1.034318980332658,1.040625533420604,1.031853063021243,1.047499474445901)
    3 = @(1.02,1.030200869198805,1.038657014648208,1.029914439594854,1.045740564816936,1.033981050292772,1.03498963189556,1.041256623974852,1.032537389848685,1.048321571600922)
    4 = @(1.02,1.030867524096433,1.039188047687956,1.030482019661197,1.046395116288663,1.034116531359075,1.035423105956784,1.041663976785046,1.03298001932781,1.048853073953561)
    5 = @(1.02,1.031147805300519,1.039411207336348,1.030720739350375,1.046670327040172,1.034173029308594,1.035605222486324,1.041834986658392,1.033166058136728,1.049076409291361)
    6 = @(1.02,1.031194866886332,1.039448671717826,1.030760827808358,1.046716538304767,1.034182488674353,1.035635793838336,1.041863685802562,1.03319729235574,1.04911390191054)
    7 = @(1.02,1.030871269151223,1.039191029899623,1.030485209017904,1.046398793522437,1.034117288088831,1.03542553986306,1.041666262776122,1.032982505355667,1.048856058597454)
    8 = @(1.02,1.029518849257446,1.038113396163743,1.029334094821863,1.045070990721332,1.033840872184959,1.034545729746843,1.040839021124068,1.032084370533425,1.047777398709599)
    9 = @(1.02,1.027136346468486,1.03621175980663,1.027309228892544,1.042732476654477,1.033339167094681,1.032991718883911,1.039373646951323,1.030500420181634,1.04587324751237)
    10 = @(1.02,1.025548442063575,1.034942228615518,1.02596174302536,1.041174339981538,1.03299489093017,1.031953268537033,1.038391610744337,1.029443594723982,1.044601552548474)
    11 = @(1.02,1.024860964033242,1.034392096100502,1.025378850034464,1.040499863339695,1.032843494906371,1.031503033497974,1.037965170676135,1.028985779773875,1.044050364923433)
    12 = @(1.02,1.024605618235262,1.034187690373553,1.025162424940486,1.040249364198649,1.03278691115623,1.031335709453685,1.03780659023645,1.028815696879221,1.043845549306438)
    13 = @(1.02,1.024660390137375,1.034231538891812,1.025208844897214,1.040303095681586,1.032799064328614,1.031371604947755,1.037840614489968,1.028852181590492,1.043889486564429)
    14 = @(1.02,1.024839856781968,1.034375201120709,1.025360958484085,1.040479156357121,1.032838824783726,1.031489204207223,1.037952066069669,1.02897172127925,1.044033436429954)
    15 = @(1.02,1.024950433976318,1.034463707963043,1.025454692273032,1.040587637339214,1.032863276339817,1.031561649512382,1.038020711014845,1.029045369658178,1.044122118106277)
    16 = @(1.02,1.025594069708065,1.034978730366155,1.026000439926246,1.041219107170272,1.033004889690496,1.031983136977711,1.038419886673475,1.029473974160547,1.044638121814537)
    17 = @(1.02,1.025997830757015,1.035301679098511,1.026342927980668,1.041615267197934,1.033093098544149,1.032247370086398,1.038669954885548,1.029742772491174,1.044961654297389)
    18 = @(1.02,1.026233346925997,1.035490009339718,1.026542751128646,1.041846360565334,1.033144325196647,1.032401436884477,1.038815698480284,1.029899538440151,1.045150313642158)
    19 = @(1.02,1.026313653414343,1.035554218231124,1.026610895081257,1.041925160797833,1.033161754142719,1.032453960189375,1.038865373426625,1.029952988246193,1.045214632824443)
    20 = @(1.02,1.025954510068854,1.035267033918498,1.026306176492177,1.041572760911851,1.033083657743922,1.03221902615455,1.038643137027541,1.029713934995425,1.04492694766599)
    21 = @(1.02,1.024787007938586,1.034332897881728,1.025316162409043,1.040427310002882,1.032827125938902,1.031454576547363,1.037919251366959,1.028936520646414,1.043991048987114)
    22 = @(1.02,1.024053034828347,1.033745210990956,1.024694206732743,1.03970730210371,1.032663817919297,1.030973435247005,1.037463065200833,1.028447555791217,1.043402149746786)
    23 = @(1.02,1.024442120149381,1.034056788695649,1.025023869097192,1.040088974464161,1.03275058159461,1.031228544786316,1.037704997542568,1.028706781637327,1.04371438008193)
    24 = @(1.02,1.025974084794576,1.035282688708521,1.026322782735872,1.041591967611223,1.033087924327514,1.032231833729968,1.038655255224555,1.029726965478879,1.044942630258398)
    25 = @(1.02,1.027752204946176,1.036703694417065,1.027832280927727,1.043336887402371,1.033470600226364,1.033393900471736,1.039753386351721,1.03091006284299,1.046365917262672)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value2 = $vals[$i]
    }
}
